# Update hotel reviews data: fill in the previously-blank
# English_Reviews_num (G2), Local_Rank (H2) and Total_Reviews_num (I2)
# cells on the hotel_info sheet.
#
# These values ("2", "466", "3") look numeric but must be stored as text
# (matching the column's existing shared-string/text cell type), so the
# number format is temporarily switched to Text ("@") before assignment
# and the cell style is restored afterwards to avoid leaving a stray
# "number stored as text" format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$rng = $ws.Range("G2:I2")
$rng.NumberFormat = "@"

$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "466"
$ws.Range("I2").Value = "3"

$rng.Style = "Normal"
